$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate the US region headers to German, keeping the other header labels intact.
$ws.Range("C1").Value = "Berg"
$ws.Range("D1").Value = "Nordosten"
$ws.Range("E1").Value = "Süd"
$ws.Range("G1").Value = "Westen"
